$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 21, pushing existing rows 21+ down by one.
$ws.Rows.Item(21).Insert()

# Set the new row's height.
$ws.Rows.Item(21).RowHeight = 60

# Fill in the new task row (row 21).
$ws.Range("A21").Value = "Task 11"
$ws.Range("B21").Value = "Sync to teacher's mode when joining the class"
$ws.Range("C21").Value = "Lian"
$ws.Range("E21").Value = "Now we set the exercise mode, but this not desirable since we don't know which mode to switch to if a student is late"
$ws.Range("F21").Value = Get-Date -Year 2022 -Month 2 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("G21").Value = 0

# Update the view state to match the saved selection.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("E21").Select()
